# Updated cryptos list with GitHub Actions - price/volume refresh
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper pattern: force-text write so numeric-looking price strings (e.g. trailing
# zeros, thousand-dot-grouped values) keep their exact literal text instead of being
# auto-coerced to a Number by the Value setter. NumberFormat is reset immediately
# after the write so no residual formatting is left on the cell.
function Set-TextValue($range, $text) {
    $range.NumberFormat = "@"
    $range.Value = $text
    $range.ClearFormats()
}

Set-TextValue $ws.Range("D2") "27.604.45"
$ws.Range("E2").Value = "  +2.91%  "
Set-TextValue $ws.Range("D3") "1.849.65"
$ws.Range("E3").Value = "  +2.02%  "
Set-TextValue $ws.Range("D4") "1.027"
$ws.Range("E4").Value = "  +2.21%  "
Set-TextValue $ws.Range("D5") "321.43"
$ws.Range("E5").Value = "  +3.89%  "
Set-TextValue $ws.Range("D6") "1.028"
$ws.Range("E6").Value = "  +2.32%  "
Set-TextValue $ws.Range("D7") "0.4371"
$ws.Range("E7").Value = "  +0.93%  "
Set-TextValue $ws.Range("D8") "0.3749"
$ws.Range("E8").Value = "  +1.11%  "
Set-TextValue $ws.Range("D9") "0.07406"
$ws.Range("E9").Value = "  +2.10%  "
Set-TextValue $ws.Range("D10") "0.8759"
$ws.Range("E10").Value = "  +1.05%  "
Set-TextValue $ws.Range("D11") "21.48"
$ws.Range("E11").Value = "  +2.56%  "
Set-TextValue $ws.Range("D12") "1.868.24"
$ws.Range("E12").Value = "  -4.19%  "
Set-TextValue $ws.Range("D13") "5.514"
$ws.Range("E13").Value = "  +2.74%  "
Set-TextValue $ws.Range("D14") "6.697"
$ws.Range("E14").Value = "  +0.32%  "
Set-TextValue $ws.Range("D15") "0.07177"
$ws.Range("E15").Value = "  +3.76%  "
Set-TextValue $ws.Range("D16") "82.69"
$ws.Range("E16").Value = "  +2.67%  "
Set-TextValue $ws.Range("D17") "1.032"
$ws.Range("E17").Value = "  +2.31%  "
Set-TextValue $ws.Range("D18") "0.000009048"
$ws.Range("E18").Value = "  +1.27%  "
Set-TextValue $ws.Range("D19") "1.027"
$ws.Range("E19").Value = "  +2.38%  "
Set-TextValue $ws.Range("D20") "15.42"
$ws.Range("E20").Value = "  +1.26%  "
Set-TextValue $ws.Range("D21") "27.629.42"
$ws.Range("E21").Value = "  +2.85%  "
Set-TextValue $ws.Range("D22") "5.263"
$ws.Range("E22").Value = "  +0.89%  "
Set-TextValue $ws.Range("D23") "11.23"
$ws.Range("E23").Value = "  +0.38%  "
Set-TextValue $ws.Range("D24") "2.078.74"
$ws.Range("E24").Value = "  -4.41%  "
$ws.Range("B25").Value = "Monero"
$ws.Range("C25").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
Set-TextValue $ws.Range("D25") "157.40"
$ws.Range("E25").Value = "  +2.32%  "
$ws.Range("B26").Value = "Toncoin"
$ws.Range("C26").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
Set-TextValue $ws.Range("D26") "1.940"
$ws.Range("E26").Value = "  +3.63%  "
Set-TextValue $ws.Range("D27") "18.74"
$ws.Range("E27").Value = "  +2.40%  "
Set-TextValue $ws.Range("D28") "5.296"
$ws.Range("E28").Value = "  +1.43%  "
Set-TextValue $ws.Range("D29") "1.938"
$ws.Range("E29").Value = "  +1.47%  "
Set-TextValue $ws.Range("D30") "116.13"
$ws.Range("E30").Value = "  +0.66%  "
Set-TextValue $ws.Range("D31") "0.09072"
$ws.Range("E31").Value = "  +1.60%  "
$ws.Range("E32").Value = "  +3.17%  "
Set-TextValue $ws.Range("D33") "0.7690"
$ws.Range("E33").Value = "  +1.60%  "
Set-TextValue $ws.Range("D34") "4.518"
$ws.Range("E34").Value = "  +1.89%  "
Set-TextValue $ws.Range("D35") "2.877"
$ws.Range("E35").Value = "  +2.57%  "
Set-TextValue $ws.Range("D36") "1.029"
$ws.Range("E36").Value = "  +2.16%  "
Set-TextValue $ws.Range("D37") "1.153"
$ws.Range("E37").Value = "  +1.85%  "
Set-TextValue $ws.Range("D38") "0.01981"
$ws.Range("E38").Value = "  +2.74%  "
Set-TextValue $ws.Range("D39") "0.05282"
$ws.Range("E39").Value = "  +1.00%  "
$ws.Range("B40").Value = "TheSandbox"
$ws.Range("C40").Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
Set-TextValue $ws.Range("D40") "0.5182"
$ws.Range("E40").Value = "  +1.94%  "
$ws.Range("B41").Value = "MXToken"
$ws.Range("C41").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
Set-TextValue $ws.Range("D41") "2.811"
$ws.Range("E41").Value = "  +4.97%  "
Set-TextValue $ws.Range("D42") "0.1674"
$ws.Range("E42").Value = "  +1.51%  "
Set-TextValue $ws.Range("D43") "6.739"
$ws.Range("E43").Value = "  +2.80%  "
Set-TextValue $ws.Range("D44") "8.606"
$ws.Range("E44").Value = "  +3.88%  "
Set-TextValue $ws.Range("D45") "108.93"
$ws.Range("E45").Value = "  +2.06%  "
$ws.Range("E46").Value = "  +2.29%  "
Set-TextValue $ws.Range("D47") "1.720"
$ws.Range("E47").Value = "  +3.97%  "
Set-TextValue $ws.Range("D48") "0.4661"
$ws.Range("E48").Value = "  +2.25%  "
$ws.Range("E49").Value = "  +1.88%  "
Set-TextValue $ws.Range("D50") "1.886"
$ws.Range("E50").Value = "  +4.30%  "
Set-TextValue $ws.Range("D51") "39.57"
$ws.Range("E51").Value = "  +5.37%  "
